# Apply cryptos list update (prices / 1h volume %) per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Worksheet, [string]$Address, [string]$Text)
    # Force the cell to store plain text even when the string looks numeric
    # (Excel would otherwise silently coerce "579.97" -> 579.97 as a number).
    $rng = $Worksheet.Range($Address)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.ClearFormats()
}

$ws.Range("D2").Value = "60.994.36"
$ws.Range("E2").Value = "  +0.92%  "
$ws.Range("D3").Value = "2.680.05"
$ws.Range("E3").Value = "  +2.98%  "
$ws.Range("E4").Value = "  -0.10%  "
Set-TextValue $ws "D5" "579.97"
$ws.Range("E5").Value = "  +1.22%  "
Set-TextValue $ws "D6" "145.57"
$ws.Range("E6").Value = "  +2.19%  "
Set-TextValue $ws "D7" "0.998"
$ws.Range("E7").Value = "  +0.11%  "
Set-TextValue $ws "D8" "0.600"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +1.64%  "
$ws.Range("E10").Value = "  +1.89%  "
Set-TextValue $ws "D11" "0.379"
$ws.Range("E11").Value = "  +4.00%  "
$ws.Range("E12").Value = "  +0.92%  "
$ws.Range("D13").Value = "3.131.27"
$ws.Range("E13").Value = "  +2.05%  "
Set-TextValue $ws "D14" "25.69"
$ws.Range("E14").Value = "  +10.60%  "
$ws.Range("D15").Value = "60.964.08"
$ws.Range("E15").Value = "  +0.85%  "
$ws.Range("E16").Value = "  +2.31%  "
$ws.Range("D17").Value = "2.671.58"
$ws.Range("E17").Value = "  +2.11%  "
Set-TextValue $ws "D18" "11.67"
$ws.Range("E18").Value = "  +2.90%  "
$ws.Range("E19").Value = "  +1.98%  "
Set-TextValue $ws "D20" "352.55"
$ws.Range("E20").Value = "  +1.77%  "
Set-TextValue $ws "D21" "6.96"
$ws.Range("E21").Value = "  -0.36%  "
Set-TextValue $ws "D22" "0.998"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E23").Value = "  +0.96%  "
Set-TextValue $ws "D24" "64.18"
$ws.Range("E24").Value = "  +1.46%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("E26").Value = "  +1.78%  "
Set-TextValue $ws "D27" "8.18"
$ws.Range("E27").Value = "  +5.69%  "
$ws.Range("E28").Value = "  +7.52%  "
$ws.Range("D29").Value = "0.0₃0819"
$ws.Range("E29").Value = "  +3.93%  "
$ws.Range("E30").Value = "  +6.38%  "
$ws.Range("E31").Value = "  +0.09%  "
Set-TextValue $ws "D32" "167.11"
$ws.Range("E32").Value = "  +3.59%  "
Set-TextValue $ws "D33" "19.91"
$ws.Range("E33").Value = "  +2.10%  "
$ws.Range("E34").Value = "  +9.10%  "
$ws.Range("E35").Value = "  +5.82%  "
$ws.Range("E36").Value = "  +8.37%  "
$ws.Range("E37").Value = "  +3.71%  "
Set-TextValue $ws "D38" "328.75"
$ws.Range("E38").Value = "  +11.59%  "
$ws.Range("E39").Value = "  +4.41%  "
Set-TextValue $ws "D40" "38.52"
$ws.Range("E40").Value = "  +1.81%  "
Set-TextValue $ws "D41" "0.883"
$ws.Range("E41").Value = "  +4.22%  "
Set-TextValue $ws "D42" "5.26"
$ws.Range("E42").Value = "  +6.61%  "
Set-TextValue $ws "D43" "20.62"
$ws.Range("E43").Value = "  +4.19%  "
Set-TextValue $ws "D44" "134.58"
$ws.Range("E44").Value = "  -2.21%  "
$ws.Range("E45").Value = "  +2.02%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws "D46" "0.617"
$ws.Range("E46").Value = "  +1.27%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws "D47" "0.0561"
$ws.Range("E47").Value = "  +3.02%  "
Set-TextValue $ws "D48" "20.56"
$ws.Range("E48").Value = "  +4.13%  "
$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws "D49" "1.00"
$ws.Range("E49").Value = "  +0.41%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws "D50" "0.0247"
$ws.Range("E50").Value = "  +3.19%  "
$ws.Range("D51").Value = "2.139.18"
$ws.Range("E51").Value = "  +5.64%  "
